$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Abril de 2020 a las 23:20"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 240511
$ws.Cells.Item(4, 3).Value = 25508
$ws.Cells.Item(4, 4).Value = 10365
$ws.Cells.Item(4, 5).Value = 224336
$ws.Cells.Item(4, 6).Value = 5421
$ws.Cells.Item(4, 7).Value = 708
$ws.Cells.Item(4, 8).Value = 5810

# Row 6
$ws.Cells.Item(6, 1).Value = "España"
$ws.Cells.Item(6, 2).Value = 112065
$ws.Cells.Item(6, 3).Value = 7947
$ws.Cells.Item(6, 4).Value = 26743
$ws.Cells.Item(6, 5).Value = 74974
$ws.Cells.Item(6, 6).Value = 6092
$ws.Cells.Item(6, 7).Value = 961
$ws.Cells.Item(6, 8).Value = 10348

# Row 7
$ws.Cells.Item(7, 1).Value = "Alemania"
$ws.Cells.Item(7, 2).Value = 84788
$ws.Cells.Item(7, 3).Value = 6807
$ws.Cells.Item(7, 4).Value = 22440
$ws.Cells.Item(7, 5).Value = 61241
$ws.Cells.Item(7, 6).Value = 3936
$ws.Cells.Item(7, 7).Value = 176
$ws.Cells.Item(7, 8).Value = 1107

# Row 9
$ws.Cells.Item(9, 1).Value = "Francia"
$ws.Cells.Item(9, 2).Value = 59105
$ws.Cells.Item(9, 3).Value = 2116
$ws.Cells.Item(9, 4).Value = 12428
$ws.Cells.Item(9, 5).Value = 41290
$ws.Cells.Item(9, 6).Value = 6399
$ws.Cells.Item(9, 7).Value = 1355
$ws.Cells.Item(9, 8).Value = 5387

# Row 17
$ws.Cells.Item(17, 1).Value = "Austria"
$ws.Cells.Item(17, 2).Value = 11123
$ws.Cells.Item(17, 3).Value = 412
$ws.Cells.Item(17, 4).Value = 1749
$ws.Cells.Item(17, 5).Value = 9216
$ws.Cells.Item(17, 6).Value = 227
$ws.Cells.Item(17, 7).Value = 12
$ws.Cells.Item(17, 8).Value = 158

# Row 23
$ws.Cells.Item(23, 1).Value = "Noruega"
$ws.Cells.Item(23, 2).Value = 5142
$ws.Cells.Item(23, 3).Value = 265
$ws.Cells.Item(23, 4).Value = 32
$ws.Cells.Item(23, 5).Value = 5060
$ws.Cells.Item(23, 6).Value = 96
$ws.Cells.Item(23, 7).Value = 6
$ws.Cells.Item(23, 8).Value = 50

# Row 24
$ws.Cells.Item(24, 1).Value = "Australia"
$ws.Cells.Item(24, 2).Value = 5139
$ws.Cells.Item(24, 3).Value = 91
$ws.Cells.Item(24, 4).Value = 585
$ws.Cells.Item(24, 5).Value = 4527
$ws.Cells.Item(24, 6).Value = 50
$ws.Cells.Item(24, 7).Value = 4
$ws.Cells.Item(24, 8).Value = 27

# Row 30
$ws.Cells.Item(30, 1).Value = "Ecuador"
$ws.Cells.Item(30, 2).Value = 3163
$ws.Cells.Item(30, 3).Value = 405
$ws.Cells.Item(30, 4).Value = 65
$ws.Cells.Item(30, 5).Value = 2978
$ws.Cells.Item(30, 6).Value = 100
$ws.Cells.Item(30, 7).Value = 22
$ws.Cells.Item(30, 8).Value = 120

# Row 35
$ws.Cells.Item(35, 1).Value = "India"
$ws.Cells.Item(35, 2).Value = 2536
$ws.Cells.Item(35, 3).Value = 538
$ws.Cells.Item(35, 4).Value = 191
$ws.Cells.Item(35, 5).Value = 2273
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 14
$ws.Cells.Item(35, 8).Value = 72

# Row 66
$ws.Cells.Item(66, 1).Value = "Marruecos"
$ws.Cells.Item(66, 2).Value = 708
$ws.Cells.Item(66, 3).Value = 54
$ws.Cells.Item(66, 4).Value = 31
$ws.Cells.Item(66, 5).Value = 633
$ws.Cells.Item(66, 6).Value = 1
$ws.Cells.Item(66, 7).Value = 5
$ws.Cells.Item(66, 8).Value = 44

# Row 103
$ws.Cells.Item(103, 1).Value = "Costa de Marfil"
$ws.Cells.Item(103, 2).Value = 194
$ws.Cells.Item(103, 3).Value = 4
$ws.Cells.Item(103, 4).Value = 15
$ws.Cells.Item(103, 5).Value = 178
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 1

# Row 146
$ws.Cells.Item(146, 1).Value = "Bermudas"
$ws.Cells.Item(146, 2).Value = 35
$ws.Cells.Item(146, 3).Value = 3
$ws.Cells.Item(146, 4).Value = 11
$ws.Cells.Item(146, 5).Value = 24
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 0

# Row 147
$ws.Cells.Item(147, 1).Value = "Guam"
$ws.Cells.Item(147, 2).Value = 32
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 0
$ws.Cells.Item(147, 5).Value = 31
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 1

# Row 167
$ws.Cells.Item(167, 1).Value = "Santa Lucia"
$ws.Cells.Item(167, 2).Value = 13
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 1
$ws.Cells.Item(167, 5).Value = 12
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 0

# Row 168
$ws.Cells.Item(168, 1).Value = "Benin"
$ws.Cells.Item(168, 2).Value = 13
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 1
$ws.Cells.Item(168, 5).Value = 12
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 0

# Row 170
$ws.Cells.Item(170, 1).Value = "Libia"
$ws.Cells.Item(170, 2).Value = 11
$ws.Cells.Item(170, 3).Value = 1
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 10
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 1
$ws.Cells.Item(170, 8).Value = 1

# Row 171
$ws.Cells.Item(171, 1).Value = "Curazao"
$ws.Cells.Item(171, 2).Value = 11
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 3
$ws.Cells.Item(171, 5).Value = 7
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 1

# Row 172
$ws.Cells.Item(172, 1).Value = "Mozambique"
$ws.Cells.Item(172, 2).Value = 10
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 10
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

# Row 173
$ws.Cells.Item(173, 1).Value = "Seychelles"
$ws.Cells.Item(173, 2).Value = 10
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 10
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

# Row 187
$ws.Cells.Item(187, 1).Value = "Santa Sede"
$ws.Cells.Item(187, 2).Value = 7
$ws.Cells.Item(187, 3).Value = 1
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 7
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

# Row 188
$ws.Cells.Item(188, 1).Value = "Fiyi"
$ws.Cells.Item(188, 2).Value = 7
$ws.Cells.Item(188, 3).Value = 2
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 7
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

